# Add 2022-Q3 data
#
# 1) Insert a new worksheet named "2022-Q3" right before the existing
#    "2022-Q2" sheet (i.e. right after "总计"), populated with the
#    quarterly fund-holding breakdown.
# 2) Update the "总计" (summary) sheet: add a new top data row for
#    2022-Q3 and keep the existing quarters below it (shifted down one
#    row), renumbering the running index in column A.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) New "2022-Q3" worksheet
# ---------------------------------------------------------------------
$anchor = $wb.Worksheets.Item("2022-Q2")
$newSheet = $wb.Worksheets.Add($anchor)
$newSheet.Name = "2022-Q3"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($col = 0; $col -lt $headers.Count; $col++) {
    $cell = $newSheet.Cells.Item(1, $col + 2)
    $cell.Value = $headers[$col]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

# code, name, scale, stock position, position ratio, market value (yi), rank
$funds = @(
    @("012421", "华夏优加生活混合A", "7.38", "90.53", "5.06", "0.3734", 1),
    @("005888", "华夏新兴消费混合A", "7.31", "89.64", "4.13", "0.3019", 2),
    @("005889", "华夏新兴消费混合C", "4.90", "89.64", "4.13", "0.2024", 2),
    @("013109", "华夏优势价值一年持有混合A", "1.87", "67.65", "3.00", "0.0561", 8),
    @("165531", "信诚多策略灵活配置混合（LOF）", "0.69", "72.58", "1.20", "0.0083", 7),
    @("012422", "华夏优加生活混合C", "0.14", "90.53", "5.06", "0.0071", 1),
    @("013110", "华夏优势价值一年持有混合C", "0.21", "67.65", "3.00", "0.0063", 8)
)

for ($i = 0; $i -lt $funds.Count; $i++) {
    $r = $i + 2
    $row = $funds[$i]

    $idxCell = $newSheet.Cells.Item($r, 1)
    $idxCell.Value = $i
    $idxCell.Font.Bold = $true
    $idxCell.HorizontalAlignment = -4108
    $idxCell.VerticalAlignment = -4160
    $idxCell.Borders.LineStyle = 1

    $newSheet.Cells.Item($r, 2).NumberFormat = "@"
    $newSheet.Cells.Item($r, 2).Value = $row[0]

    $newSheet.Cells.Item($r, 3).NumberFormat = "@"
    $newSheet.Cells.Item($r, 3).Value = $row[1]

    $newSheet.Cells.Item($r, 4).NumberFormat = "@"
    $newSheet.Cells.Item($r, 4).Value = $row[2]

    $newSheet.Cells.Item($r, 5).NumberFormat = "@"
    $newSheet.Cells.Item($r, 5).Value = $row[3]

    $newSheet.Cells.Item($r, 6).NumberFormat = "@"
    $newSheet.Cells.Item($r, 6).Value = $row[4]

    $newSheet.Cells.Item($r, 7).NumberFormat = "@"
    $newSheet.Cells.Item($r, 7).Value = $row[5]

    $newSheet.Cells.Item($r, 8).Value = $row[6]
}

# ---------------------------------------------------------------------
# 2) "总计" summary sheet - prepend the 2022-Q3 row, shift the rest down
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

# date label, held fund count, held market value (yi)
$summaryRows = @(
    @("2022-Q3", 7, 0.96),
    @("2022-Q2", 8, 1.2),
    @("2022-Q1", 10, 3.12),
    @("2021-Q4", 3, 2.01),
    @("2021-Q3", 4, 0.94),
    @("2021-Q2", 10, 0.7),
    @("2021-Q1", 22, 1.23),
    @("2020-Q4", 2, 0.57)
)

for ($i = 0; $i -lt $summaryRows.Count; $i++) {
    $r = $i + 2
    $row = $summaryRows[$i]

    $idxCell = $totalSheet.Cells.Item($r, 1)
    $idxCell.Value = $i
    $idxCell.Font.Bold = $true
    $idxCell.HorizontalAlignment = -4108
    $idxCell.VerticalAlignment = -4160
    $idxCell.Borders.LineStyle = 1

    $totalSheet.Cells.Item($r, 2).Value = $row[0]
    $totalSheet.Cells.Item($r, 3).Value = $row[1]
    $totalSheet.Cells.Item($r, 4).Value = $row[2]
}

Write-Output "2022-Q3 sheet added; 总计 summary updated."
